# Updated legacy GSC export data:
#  - Drop the oldest date row (2025-10-12) and shift every remaining
#    date/row up by one position.
#  - Append a brand-new trailing date (2026-01-10) with a fresh (zero)
#    HTTPS-URL count, continuing the existing (and now shifted) series.
#
# The "Chart" sheet holds the Date / Non-HTTPS URLs / HTTPS URLs table
# in rows 2..91 (row 1 is the header). The "Table" sheet only has
# headers and is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstDataRow = 2
$lastDataRow = 91

# Shift rows (firstDataRow+1)..lastDataRow up into firstDataRow..(lastDataRow-1).
# Using Range.Copy (rather than re-typing the text) preserves the cells'
# existing string typing so the dates are not reinterpreted/auto-converted
# into Excel date serial numbers.
$srcRange = $ws.Range("A" + ($firstDataRow + 1) + ":C" + $lastDataRow)
$dstRange = $ws.Range("A" + $firstDataRow + ":C" + ($lastDataRow - 1))
$srcRange.Copy($dstRange)

# Build the new trailing date as literal text in a scratch cell (so it is
# stored as a shared string, just like the other date cells) and copy it
# into place, then set its HTTPS URL count.
$scratch = $ws.Cells.Item($lastDataRow + 100, 10)
$scratch.Formula = '="2026-01-10"'
$scratch.Copy($ws.Cells.Item($lastDataRow, 1))
$scratch.ClearContents()

$ws.Cells.Item($lastDataRow, 3).Value = 0
